$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before the "Late" column (column N), shifting the
# existing "Late" / "Date" / "Outstanding" columns one to the right.
$ws.Columns("N").Insert()

# The newly inserted column takes on the width of its left neighbour ("In Advance").
$ws.Columns("N").ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet (it was "Transactions" before) and
# move the selection to R7, matching the new layout.
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
